{"js": "// Update the division-fact worksheet table: replace each populated cell's\n// text with the new value, in document order, while preserving existing\n// run/paragraph formatting (font, size, alignment).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index -> ordered [oldText, newText] pairs for the 5 cells in that row.\nconst rowEdits = {\n  0: [\"66\u00f76=\", \"73\u00f72=\", \"76\u00f75=\", \"33\u00f73=\", \"42\u00f74=\"],\n  4: [\"42\u00f75=\", \"25\u00f74=\", \"38\u00f75=\", \"92\u00f74=\", \"41\u00f79=\"],\n  8: [\"60\u00f76=\", \"36\u00f73=\", \"53\u00f78=\", \"96\u00f72=\", \"77\u00f79=\"],\n  12: [\"71\u00f73=\", \"97\u00f78=\", \"45\u00f78=\", \"92\u00f74=\", \"38\u00f78=\"],\n  16: [\"40\u00f73=\", \"98\u00f74=\", \"97\u00f74=\", \"17\u00f72=\", \"24\u00f76=\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowEdits)) {\n  const rowIndex = Number(rowIndexStr);\n  const newTexts = rowEdits[rowIndex];\n  for (let col = 0; col < newTexts.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const range = cell.body.getRange();\n    range.insertText(newTexts[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-fact worksheet table: replace each populated cell's\n# text with the new value, row by row (Word COM Table.Cell is 1-indexed),\n# while preserving existing run/paragraph formatting (font, size, alignment).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowEdits = @{\n    1  = @(\"66\u00f76=\", \"73\u00f72=\", \"76\u00f75=\", \"33\u00f73=\", \"42\u00f74=\")\n    5  = @(\"42\u00f75=\", \"25\u00f74=\", \"38\u00f75=\", \"92\u00f74=\", \"41\u00f79=\")\n    9  = @(\"60\u00f76=\", \"36\u00f73=\", \"53\u00f78=\", \"96\u00f72=\", \"77\u00f79=\")\n    13 = @(\"71\u00f73=\", \"97\u00f78=\", \"45\u00f78=\", \"92\u00f74=\", \"38\u00f78=\")\n    17 = @(\"40\u00f73=\", \"98\u00f74=\", \"97\u00f74=\", \"17\u00f72=\", \"24\u00f76=\")\n}\n\nforeach ($rowIndex in $rowEdits.Keys) {\n    $newTexts = $rowEdits[$rowIndex]\n    for ($col = 1; $col -le $newTexts.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $newTexts[$col - 1]\n    }\n}\n"}
